$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows 2-21 down to 3-22)
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the new age group SSPM0_4 / SSPF0_4
$ws.Cells.Item(2, 1).Value = "SSPM0_4"
$ws.Cells.Item(2, 2).Value = "Chil4_8"
$ws.Cells.Item(2, 3).Value = "SSPF0_4"
$ws.Cells.Item(2, 4).Value = "Chil4_8"

# Update the active selection to C2 to match the author's final state
$ws.Range("C2").Select()
